$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) cell format used by all
# B/C/D/E data cells in this sheet. Used to restore the style after
# temporarily forcing a text number-format, so that plain-number-looking
# strings (e.g. "213.99") are not silently converted to real numbers by
# Excel, while keeping the original styling (no explicit "s" attribute).
$defaultStyle = $ws.Cells.Item(4, 4).Style

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $defaultStyle
}

Set-TextCell 2 4 '25.961.64'
Set-TextCell 2 5 '  -0.15%  '

Set-TextCell 3 4 '1.627.49'
Set-TextCell 3 5 '  -0.90%  '

Set-TextCell 4 5 '  +0.03%  '

Set-TextCell 5 4 '213.99'
Set-TextCell 5 5 '  -0.84%  '

Set-TextCell 6 5 '  -0.90%  '

Set-TextCell 7 5 '  +0.05%  '

Set-TextCell 8 4 '0.251'
Set-TextCell 8 5 '  -1.86%  '

Set-TextCell 9 4 '0.0618'
Set-TextCell 9 5 '  -3.27%  '

Set-TextCell 10 4 '18.38'
Set-TextCell 10 5 '  -5.98%  '

Set-TextCell 11 4 '0.0789'
Set-TextCell 11 5 '  -0.87%  '

Set-TextCell 12 2 'WrappedliquidstakedEther2.0'
Set-TextCell 12 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 12 4 '1.853.53'
Set-TextCell 12 5 '  -0.89%  '

Set-TextCell 13 2 'WrappedEther'
Set-TextCell 13 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 13 4 '1.635.99'
Set-TextCell 13 5 '  -0.86%  '

Set-TextCell 14 5 '  -2.11%  '

Set-TextCell 15 4 '0.526'
Set-TextCell 15 5 '  -3.31%  '

Set-TextCell 16 4 '25.948.41'
Set-TextCell 16 5 '  -0.30%  '

Set-TextCell 17 4 '0.0₃0738'
Set-TextCell 17 5 '  -3.35%  '

Set-TextCell 18 4 '61.27'
Set-TextCell 18 5 '  -3.43%  '

Set-TextCell 19 5 '  +0.08%  '

Set-TextCell 20 4 '192.10'
Set-TextCell 20 5 '  -1.09%  '

Set-TextCell 21 5 '  -2.89%  '

Set-TextCell 22 4 '9.58'
Set-TextCell 22 5 '  -3.50%  '

Set-TextCell 23 4 '6.07'
Set-TextCell 23 5 '  -2.20%  '

Set-TextCell 24 4 '0.133'
Set-TextCell 24 5 '  +0.35%  '

Set-TextCell 25 4 '143.64'
Set-TextCell 25 5 '  +0.43%  '

Set-TextCell 26 5 '  +0.08%  '

Set-TextCell 27 5 '  -3.71%  '

Set-TextCell 28 4 '6.72'
Set-TextCell 28 5 '  -2.30%  '

Set-TextCell 29 4 '15.19'
Set-TextCell 29 5 '  -2.11%  '

Set-TextCell 30 5 '  -1.09%  '

Set-TextCell 31 5 '  -2.35%  '

Set-TextCell 32 4 '3.12'
Set-TextCell 32 5 '  -4.40%  '

Set-TextCell 33 4 '3.11'
Set-TextCell 33 5 '  -5.61%  '

Set-TextCell 34 5 '  -2.53%  '

Set-TextCell 35 5 '  -2.77%  '

Set-TextCell 36 4 '1.125.57'
Set-TextCell 36 5 '  -0.50%  '

Set-TextCell 37 4 '0.849'
Set-TextCell 37 5 '  -6.18%  '

Set-TextCell 38 5 '  -1.19%  '

Set-TextCell 39 4 '0.520'
Set-TextCell 39 5 '  -3.84%  '

Set-TextCell 40 5 '  -2.44%  '

Set-TextCell 41 4 '98.11'
Set-TextCell 41 5 '  -0.97%  '

Set-TextCell 42 2 'RocketPoolETH'
Set-TextCell 42 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 42 4 '1.762.73'
Set-TextCell 42 5 '  -0.94%  '

Set-TextCell 43 2 'TrustWalletToken'
Set-TextCell 43 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 43 4 '0.763'
Set-TextCell 43 5 '  -4.37%  '

Set-TextCell 44 4 '5.17'
Set-TextCell 44 5 '  -5.52%  '

Set-TextCell 45 5 '  -9.95%  '

Set-TextCell 46 4 '0.0531'
Set-TextCell 46 5 '  +1.83%  '

Set-TextCell 47 4 '54.32'
Set-TextCell 47 5 '  -3.91%  '

Set-TextCell 48 5 '  -1.27%  '

Set-TextCell 49 5 '  -0.29%  '

Set-TextCell 50 4 '1.00'
Set-TextCell 50 5 '  +0.14%  '

Set-TextCell 51 4 '7.45'
Set-TextCell 51 5 '  -3.94%  '
